# Update "想去人数" (F) and "最低票价" (G) figures on the 展览 and 全部类型
# sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 6631
    $ws.Range("G2").Value = 75

    $ws.Range("F6").Value = 2023

    $ws.Range("F7").Value = 1537

    $ws.Range("F9").Value = 1010

    $ws.Range("F10").Value = 427

    $ws.Range("F11").Value = 16

    $ws.Range("F12").Value = 5634
}
